$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "596.00", "1.00")
# must be forced to Text format first, otherwise Excel auto-converts them
# to numeric values and strips formatting (e.g. trailing zeros).
# NumberFormat/Style must be handled per-cell (a Union range does not
# reliably propagate to every member cell in this engine).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '596.00'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.91'
$ws.Range('D6').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.91'
$ws.Range('D13').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.33'
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '358.21'
$ws.Range('D19').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.39'
$ws.Range('D21').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '74.59'
$ws.Range('D24').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.75'
$ws.Range('D26').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '555.81'
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.95'
$ws.Range('D31').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.126'
$ws.Range('D35').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '159.59'
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.65'
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.369'
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.31'
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.59'
$ws.Range('D42').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '156.31'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.77'
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '21.88'
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.67'
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0772'
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.612'
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.566'
$ws.Range('D51').Style = "Normal"

# Cells whose new values are safely interpreted/stored as text as-is
$ws.Range('D2').Value = '68.190.96'
$ws.Range('D3').Value = '2.641.87'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('D14').Value = '3.123.93'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('E15').Value = '  -2.94%  '
$ws.Range('D16').Value = '68.179.98'
$ws.Range('D17').Value = '2.638.21'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  -3.89%  '
$ws.Range('E23').Value = '  -1.10%  '
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('D27').Value = '2.832.03'
$ws.Range('E27').Value = '  +2.07%  '
$ws.Range('E28').Value = '  -3.41%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -3.43%  '
$ws.Range('E31').Value = '  -2.44%  '
$ws.Range('E32').Value = '  -3.86%  '
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('E35').Value = '  -4.18%  '
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +1.25%  '
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('E40').Value = '  -2.99%  '
$ws.Range('E41').Value = '  -2.37%  '
$ws.Range('E42').Value = '  -3.06%  '
$ws.Range('E43').Value = '  -5.83%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('E49').Value = '  -2.24%  '
$ws.Range('E51').Value = '  -0.66%  '
